$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data stores every value (dates, percentages, numbers, units, ...) as
# plain text, not as native numeric/date/percentage cell types. Force the new rows
# to use the Text number format before writing values so the engine does not
# auto-convert number-, percentage- or date-looking strings (e.g. "42 %", "50", 
# "2022-01-17 09:45:53") into numeric/date values.
$ws.Range("A17:BI18").NumberFormat = "@"

# Row 17
$ws.Cells.Item(17, 1).Value = '2022-01-17 09:45:53'
$ws.Cells.Item(17, 2).Value = '16.9 °C'
$ws.Cells.Item(17, 3).Value = '24.5 °C'
$ws.Cells.Item(17, 4).Value = '1.3 °C'
$ws.Cells.Item(17, 5).Value = '**.* °C'
$ws.Cells.Item(17, 6).Value = '23.7 °C'
$ws.Cells.Item(17, 7).Value = '**.* °C'
$ws.Cells.Item(17, 8).Value = '42 %'
$ws.Cells.Item(17, 9).Value = '**.* %'
$ws.Cells.Item(17, 10).Value = '50 %'
$ws.Cells.Item(17, 11).Value = '50 %'
$ws.Cells.Item(17, 12).Value = '0 mV'
$ws.Cells.Item(17, 13).Value = '0 mV'
$ws.Cells.Item(17, 14).Value = '50 %'
$ws.Cells.Item(17, 15).Value = '50 %'
$ws.Cells.Item(17, 16).Value = '**.* mV'
$ws.Cells.Item(17, 17).Value = '**.* mV'
$ws.Cells.Item(17, 18).Value = '100 %'
$ws.Cells.Item(17, 19).Value = '0 %'
$ws.Cells.Item(17, 20).Value = '100 %'
$ws.Cells.Item(17, 21).Value = '0 %'
$ws.Cells.Item(17, 22).Value = '100 %'
$ws.Cells.Item(17, 23).Value = '32 %'
$ws.Cells.Item(17, 24).Value = '68 %'
$ws.Cells.Item(17, 25).Value = '#.## g/m³'
$ws.Cells.Item(17, 26).Value = '8.97 g/m³'
$ws.Cells.Item(17, 27).Value = 'NORMAL'
$ws.Cells.Item(17, 28).Value = '7:00-22:00'
$ws.Cells.Item(17, 29).Value = 'AWAY              22:00'
$ws.Cells.Item(17, 30).Value = 'NORMAL'
$ws.Cells.Item(17, 31).Value = '17.0 °C'
$ws.Cells.Item(17, 32).Value = '24.5 °C'
$ws.Cells.Item(17, 33).Value = '1.5 °C'
$ws.Cells.Item(17, 34).Value = '50'
$ws.Cells.Item(17, 35).Value = '50 %'
$ws.Cells.Item(17, 36).Value = '50 %'
$ws.Cells.Item(17, 37).Value = '50 %'
$ws.Cells.Item(17, 38).Value = '50 %'
$ws.Cells.Item(17, 39).Value = '32 %'
$ws.Cells.Item(17, 40).Value = '89 %'
$ws.Cells.Item(17, 41).Value = '1110 W'
$ws.Cells.Item(17, 42).Value = '537 W'
$ws.Cells.Item(17, 43).Value = '500 W'
$ws.Cells.Item(17, 44).Value = '0.23'
$ws.Cells.Item(17, 45).Value = '0.28'
$ws.Cells.Item(17, 46).Value = '2.14 kWh'
$ws.Cells.Item(17, 47).Value = '292.39 kWh'
$ws.Cells.Item(17, 48).Value = '2677.95 kWh'
$ws.Cells.Item(17, 49).Value = '1.90 kWh'
$ws.Cells.Item(17, 50).Value = '270.61 kWh'
$ws.Cells.Item(17, 51).Value = '2128.48 kWh'
$ws.Cells.Item(17, 52).Value = '4.88 kWh'
$ws.Cells.Item(17, 53).Value = '656.20 kWh'
$ws.Cells.Item(17, 54).Value = '8237.50 kWh'
$ws.Cells.Item(17, 55).Value = '19.0 °C'
$ws.Cells.Item(17, 56).Value = '--.- °C'
$ws.Cells.Item(17, 57).Value = '--.- %'
$ws.Cells.Item(17, 58).Value = '42 %'
$ws.Cells.Item(17, 59).Value = '--.- %'
$ws.Cells.Item(17, 60).Value = '42 %'
$ws.Cells.Item(17, 61).Value = '144982044'

# Row 18
$ws.Cells.Item(18, 1).Value = '2022-01-18 15:10:41'
$ws.Cells.Item(18, 2).Value = '13.9 °C'
$ws.Cells.Item(18, 3).Value = '23.8 °C'
$ws.Cells.Item(18, 4).Value = '-1.5 °C'
$ws.Cells.Item(18, 5).Value = '**.* °C'
$ws.Cells.Item(18, 6).Value = '24.2 °C'
$ws.Cells.Item(18, 7).Value = '**.* °C'
$ws.Cells.Item(18, 8).Value = '37 %'
$ws.Cells.Item(18, 9).Value = '**.* %'
$ws.Cells.Item(18, 10).Value = '50 %'
$ws.Cells.Item(18, 11).Value = '50 %'
$ws.Cells.Item(18, 12).Value = '0 mV'
$ws.Cells.Item(18, 13).Value = '0 mV'
$ws.Cells.Item(18, 14).Value = '50 %'
$ws.Cells.Item(18, 15).Value = '50 %'
$ws.Cells.Item(18, 16).Value = '**.* mV'
$ws.Cells.Item(18, 17).Value = '**.* mV'
$ws.Cells.Item(18, 18).Value = '100 %'
$ws.Cells.Item(18, 19).Value = '0 %'
$ws.Cells.Item(18, 20).Value = '100 %'
$ws.Cells.Item(18, 21).Value = '0 %'
$ws.Cells.Item(18, 22).Value = '100 %'
$ws.Cells.Item(18, 23).Value = '32 %'
$ws.Cells.Item(18, 24).Value = '71 %'
$ws.Cells.Item(18, 25).Value = '#.## g/m³'
$ws.Cells.Item(18, 26).Value = '8.14 g/m³'
$ws.Cells.Item(18, 27).Value = 'NORMAL'
$ws.Cells.Item(18, 28).Value = '7:00-22:00'
$ws.Cells.Item(18, 29).Value = 'AWAY              22:00'
$ws.Cells.Item(18, 30).Value = 'NORMAL'
$ws.Cells.Item(18, 31).Value = '13.9 °C'
$ws.Cells.Item(18, 32).Value = '23.8 °C'
$ws.Cells.Item(18, 33).Value = '-1.7 °C'
$ws.Cells.Item(18, 34).Value = '50'
$ws.Cells.Item(18, 35).Value = '50 %'
$ws.Cells.Item(18, 36).Value = '50 %'
$ws.Cells.Item(18, 37).Value = '50 %'
$ws.Cells.Item(18, 38).Value = '50 %'
$ws.Cells.Item(18, 39).Value = '32 %'
$ws.Cells.Item(18, 40).Value = '90 %'
$ws.Cells.Item(18, 41).Value = '1233 W'
$ws.Cells.Item(18, 42).Value = '536 W'
$ws.Cells.Item(18, 43).Value = '500 W'
$ws.Cells.Item(18, 44).Value = '0.22'
$ws.Cells.Item(18, 45).Value = '0.26'
$ws.Cells.Item(18, 46).Value = '5.15 kWh'
$ws.Cells.Item(18, 47).Value = '297.83 kWh'
$ws.Cells.Item(18, 48).Value = '2690.30 kWh'
$ws.Cells.Item(18, 49).Value = '4.71 kWh'
$ws.Cells.Item(18, 50).Value = '275.86 kWh'
$ws.Cells.Item(18, 51).Value = '2139.89 kWh'
$ws.Cells.Item(18, 52).Value = '11.89 kWh'
$ws.Cells.Item(18, 53).Value = '663.67 kWh'
$ws.Cells.Item(18, 54).Value = '8264.42 kWh'
$ws.Cells.Item(18, 55).Value = '19.0 °C'
$ws.Cells.Item(18, 56).Value = '--.- °C'
$ws.Cells.Item(18, 57).Value = '--.- %'
$ws.Cells.Item(18, 58).Value = '37 %'
$ws.Cells.Item(18, 59).Value = '--.- %'
$ws.Cells.Item(18, 60).Value = '37 %'
$ws.Cells.Item(18, 61).Value = '144982044'

